# 26 Apr 2024 - Codes were modified to support one group of questions in question bank.
# Add two new candidate rows (VM trade, batch 2222) to the candidates sheet
# and mark the batch_no column (B) with the "General" number format used
# on the rest of that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 12 ---
$ws.Range("A12").Value = "VM"
$ws.Range("B12").Value = 2222
$ws.Range("C12").Value = "VM98765"
$ws.Range("D12").Value = "Wong"
$ws.Range("E12").Value = "Tai Sin"
$ws.Range("F12").Value = 65432109
$ws.Range("F12").NumberFormat = "@"

# --- New row 13 ---
$ws.Range("A13").Value = "VM"
$ws.Range("B13").Value = 2222
$ws.Range("C13").Value = "VM54321"
$ws.Range("D13").Value = "Chan"
$ws.Range("E13").Value = "Tai Man"
$ws.Range("F13").Value = 98765432
$ws.Range("F13").NumberFormat = "@"

# --- Apply the (newly introduced) number format to the whole batch_no column ---
$ws.Range("B2:B13").NumberFormat = "General"

# --- Update selection to match the authored workbook ---
$ws.Range("B2:B13").Select() | Out-Null
